$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.465.30'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '2.927.16'
$ws.Range('E3').Value = '  +4.71%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''352.68'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = '''113.02'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '''40.19'
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  +3.73%  '
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '''20.18'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').Value = '3.384.53'
$ws.Range('E15').Value = '  +4.84%  '
$ws.Range('D16').Value = '2.941.47'
$ws.Range('E16').Value = '  +5.48%  '
$ws.Range('D17').Value = '''0.994'
$ws.Range('E17').Value = '  +5.35%  '
$ws.Range('D18').Value = '52.450.05'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''7.73'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').Value = '''3.34'
$ws.Range('E20').Value = '  +5.03%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = '''14.46'
$ws.Range('E21').Value = '  +7.01%  '
$ws.Range('D22').Value = '0.0₃0984'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').Value = '''71.19'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').Value = '''271.57'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('E25').Value = '  +2.83%  '
$ws.Range('D26').Value = '''26.97'
$ws.Range('E26').Value = '  +3.62%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  +3.15%  '
$ws.Range('D30').Value = '''38.18'
$ws.Range('E30').Value = '  +3.83%  '
$ws.Range('D31').Value = '''6.54'
$ws.Range('E31').Value = '  +5.90%  '
$ws.Range('E32').Value = '  +14.05%  '
$ws.Range('D33').Value = '''6.23'
$ws.Range('E33').Value = '  +9.16%  '
$ws.Range('D34').Value = '''53.33'
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').Value = '''0.0940'
$ws.Range('E35').Value = '  +10.13%  '
$ws.Range('D36').Value = '''0.0453'
$ws.Range('E36').Value = '  +2.94%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '''3.35'
$ws.Range('E38').Value = '  +7.49%  '
$ws.Range('E39').Value = '  +4.81%  '
$ws.Range('D40').Value = '''18.80'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').Value = '''2.73'
$ws.Range('E41').Value = '  +10.47%  '
$ws.Range('D42').Value = '''24.51'
$ws.Range('E42').Value = '  +13.00%  '
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('D44').Value = '''123.02'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''2.60'
$ws.Range('E46').Value = '  +8.28%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.217.93'
$ws.Range('E47').Value = '  +4.42%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''3.56'
$ws.Range('E48').Value = '  +6.03%  '
$ws.Range('D49').Value = '''0.263'
$ws.Range('E49').Value = '  +24.51%  '
$ws.Range('D50').Value = '''0.0338'
$ws.Range('E50').Value = '  +16.17%  '
$ws.Range('D51').Value = '''0.962'
$ws.Range('E51').Value = '  +5.46%  '
